$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping the cell as plain text
# (matches the source data, which stores all Price/rank-string values as text).
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "245.96"
Set-TextValue "D4" "5.434"
Set-TextValue "D5" "0.05839"
Set-TextValue "D6" "3.374"
Set-TextValue "D7" "6.347"
Set-TextValue "D8" "0.8088"
Set-TextValue "D9" "0.9687"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01124"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1428"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07457"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03456"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03039"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "4.166"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09401"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001604"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04822"
$ws.Range("E18").Value = "17CoinExTokenCET"
Set-TextValue "D19" "0.006097"
Set-TextValue "D21" "0.0009975"
Set-TextValue "D23" "3.703"
Set-TextValue "D24" "2.214"
Set-TextValue "D25" "0.3207"
Set-TextValue "D40" "0.03860"
Set-TextValue "D41" "0.006732"
Set-TextValue "D42" "0.1077"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.006682"
Set-TextValue "D45" "0.00005607"
Set-TextValue "D48" "0.1457"
Set-TextValue "D49" "0.00002101"
